# Apply updated accuracy figures from the new 201006 run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: logistic_embeddings
$ws.Range("C5").Value = 0.529
$ws.Range("D5").Value = 0.612
$ws.Range("E5").Value = 0.631
$ws.Range("F5").Value = 0.708
$ws.Range("G5").Value = 0.705
$ws.Range("H5").Value = 0.718

# Row 7: classical-best-embeddings -> classical-best-embed
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.529
$ws.Range("F7").Value = 0.708
$ws.Range("G7").Value = 0.705
$ws.Range("H7").Value = 0.718

# Row 8: BERT-base
$ws.Range("C8").Value = 0.515
$ws.Range("E8").Value = 0.696
$ws.Range("F8").Value = 0.716
$ws.Range("G8").Value = 0.742
$ws.Range("H8").Value = 0.757

# Row 9: BERT-base-nli
$ws.Range("B9").Value = 0.266
$ws.Range("C9").Value = 0.472
$ws.Range("D9").Value = 0.624
$ws.Range("E9").Value = 0.658
$ws.Range("F9").Value = 0.654
$ws.Range("G9").Value = 0.69
$ws.Range("H9").Value = 0.7
